$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "HORA" (time) column before the existing CPF column (C),
# shifting CPF/ENDEREÇO/CTPS one column to the right (D/E/F).
$ws.Columns("C").Insert()

# Header for the new column
$ws.Range("C1").Value = "HORA"

# Time-of-day values for the two data rows (stored as Excel day fractions)
$ws.Range("C2").Value = [DateTime]::FromOADate(0.47916666666666669)
$ws.Range("C3").Value = [DateTime]::FromOADate(0.72916666666666663)

# Format the new column as a time (h:mm)
$ws.Range("C2:C3").NumberFormat = "h:mm"

# Match the new column's width to the neighboring DATA column
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Update the active selection to reflect the edited cell
$ws.Range("C3").Select()
